# issue #5: stock data from json to db
#
# The "股票" (stock) sheet (sheet 4) gains three new columns:
#   - "category"    inserted right after "property_category" (new column I),
#                    value "normal" for every data row (matches the
#                    output/normal/... path segment of the source file)
#   - "source_file" appended after "legislator_id" (new column M),
#                    value "tmp845a1" (matches the *_tmp845a1.xlsx file stem)
#   - "index"       appended after "source_file" (new column N),
#                    value equal to the row's original identifier (col A)
#
# Existing columns keep their values; they're simply shifted right by the
# "category" insertion so the header order becomes:
#   name, owner, quantity, face_value, currency, total, property_category,
#   category, date, legislator_name, legislator_id, source_file, index

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# Insert the new "category" column between H (property_category) and the old
# I (date); this shifts the old I/J/K (date/legislator_name/legislator_id)
# one column to the right while carrying their formatting along.
$ws.Columns.Item(9).Insert()

# Append two brand-new columns for "source_file" and "index" right after the
# (now shifted) legislator_id column (L); inserting as columns (rather than
# just writing into bare cells) makes Excel copy the neighbouring column's
# formatting onto them.
$ws.Columns.Item(13).Insert()
$ws.Columns.Item(14).Insert()

$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

$lastRow = 18
For ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Value = "normal"
    $ws.Cells.Item($r, 13).Value = "tmp845a1"
    $ws.Cells.Item($r, 14).Value = $ws.Cells.Item($r, 1).Value2
}
